$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STATUS")

# Update row 2 in place: Bob -> Doni, all stats reset to zero (mock/test data).
$ws.Range("A2").Value = "Doni"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# B2/E2/F2 hold percentage-looking text ("0%") that must stay a literal
# string rather than Excel's usual auto-convert-to-number-with-% behaviour.
# Route the text through a formula result, then copy/paste-values it back
# onto itself: the pasted value keeps the formula's text type (no percent
# auto-detection) while leaving the original cell style untouched.
$ws.Range("B2").Formula = "=""0%"""
$ws.Range("E2").Formula = "=""0%"""
$ws.Range("F2").Formula = "=""0%"""
$ws.Range("B2:F2").Copy()
$ws.Range("B2:F2").PasteSpecial(-4163)

# Remove the old row 3 (previously "Doni") entirely, shrinking the used range to A1:F2
$ws.Rows.Item(3).Delete()
